# Weekly update: insert a new daily price record as row 74, pushing the
# existing rows 74..127 down to 75..128 (last row becomes 128).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 74 (shifts rows 74-127 -> 75-128).
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with the new observation.
$ws.Range("A74").Value = 5
$ws.Range("B74").Value = 'Macroferia Regional de Talca'
$ws.Range("C74").Value = 'Maule'
$ws.Range("D74").Value = 44907
$ws.Range("E74").Value = 7
$ws.Range("F74").Value = 100112022
$ws.Range("G74").Value = 'Arveja Verde'
$ws.Range("H74").Value = 'Sin especificar'
$ws.Range("I74").Value = 'Primera'
$ws.Range("J74").Value = 500
$ws.Range("K74").Value = 25000
$ws.Range("L74").Value = 25000
$ws.Range("M74").Value = 25000
$ws.Range("N74").Value = '$/saco 25 kilos'
$ws.Range("O74").Value = 'Carahue'
$ws.Range("P74").Value = 1000
$ws.Range("Q74").Value = 25
$ws.Range("R74").Value = 'Hortaliza'
